# Add "2022-Q1" worksheet (new fund-holding detail sheet), inserted right
# before the existing "总计" (total) summary sheet, and record the new
# quarter's totals into the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Locate the existing "2021-Q4" sheet (style/layout template) and the
#    "总计" sheet (insertion anchor + the one we need to update).
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalName = $totalSheetBefore.Name

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q1" sheet right before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Re-fetch the "总计" sheet by name: its positional index shifted by one
# once the new sheet was inserted in front of it.
$totalSheet = $wb.Worksheets.Item($totalName)

# Copy the header row formatting (bold + border + centered) from the
# 2021-Q4 template so the new sheet matches the look of the others.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the "index column" formatting (bold border, centered) used for
# column A on the template sheet.
$templateSheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# B2:G3 hold values that look numeric ("002379", "1.84", ...) but must be
# stored as plain text (leading zeros, fixed decimal formatting). Force
# text entry via a temporary "@" number format, then strip the format
# back off so the cells end up as plain, unstyled text cells (matching
# the rest of the data cells on this sheet).
$textCells = $newSheet.Range("B2:G3")
$textCells.NumberFormat = "@"

$newSheet.Range("B2").Value = "002379"
$newSheet.Range("C2").Value = "工银瑞信香港中小盘股票（QDII）人民币"
$newSheet.Range("D2").Value = "1.84"
$newSheet.Range("E2").Value = "86.48"
$newSheet.Range("F2").Value = "5.10"
$newSheet.Range("G2").Value = "0.0938"
$newSheet.Range("H2").Value = 2

$newSheet.Range("B3").Value = "002380"
$newSheet.Range("C3").Value = "工银瑞信香港中小盘股票（QDII）美元"
$newSheet.Range("D3").Value = "1.84"
$newSheet.Range("E3").Value = "86.48"
$newSheet.Range("F3").Value = "5.10"
$newSheet.Range("G3").Value = "0.0938"
$newSheet.Range("H3").Value = 2

$textCells.ClearFormats()

# ---------------------------------------------------------------------
# 3. Update the "总计" sheet: add a new row for 2022-Q1 at the top of the
#    data (row 2), push the previous rows down by one, and renumber the
#    index column (A).
# ---------------------------------------------------------------------

# Grab current values up-front (before anything is overwritten).
$dates = @()
$counts = @()
$values = @()
for ($r = 2; $r -le 6; $r++) {
    $dates += , $totalSheet.Cells.Item($r, 2).Value2
    $counts += , $totalSheet.Cells.Item($r, 3).Value2
    $values += , $totalSheet.Cells.Item($r, 4).Value2
}

# Extend the index-column formatting down to the new last row (row 7).
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

# Prepend the new quarter, shifting everything else down one row.
$newDates = @("2022-Q1") + $dates
$newCounts = @(2) + $counts
$newValues = @(0.19) + $values

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $newDates[$i]
    $totalSheet.Cells.Item($r, 3).Value = $newCounts[$i]
    $totalSheet.Cells.Item($r, 4).Value = $newValues[$i]
}
